# Auto-update draw results: append the latest Pick 3 draw as a new row
# at the bottom of the "Results" sheet (mirrors the site's daily export).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40

# Columns A (date) and C (phase code) look numeric/date-like, so Excel would
# otherwise auto-coerce them to a date serial / number on assignment.
# Pre-formatting as Text keeps them as literal strings, matching the rest
# of the column.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-10-26"
$ws.Cells.Item($row, 2).Value = "Pick 3"
$ws.Cells.Item($row, 3).Value = "251026"
$ws.Cells.Item($row, 4).Value = "5-9-9"
$ws.Cells.Item($row, 5).Value = "2025-10-26T21:36:29.707+04:00"
